$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - updated TPM-derived values
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 0.03329446086555556
$ws.Range("R2").Value = 0.29965014779
$ws.Range("S2").Value = 0.02073452941466921
$ws.Range("T2").Value = 0.02073452941466921

# Row 3 - updated TPM-derived values (M3, N3, Q3 unchanged)
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("R3").Value = 8.134439651099999
$ws.Range("S3").Value = 0.5628689972673966
$ws.Range("T3").Value = 0.5628689972673966

# Row 4 - updated TPM-derived values
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 0.6686284413877777
$ws.Range("R4").Value = 6.017655972489999
$ws.Range("S4").Value = 0.4163964733179342
$ws.Range("T4").Value = 0.4163964733179341
